# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The G column (header "K" in G1) is regenerated from an updated source;
# write the new per-row values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 2
    6  = 3
    7  = 2
    8  = 1
    9  = 2
    10 = 4
    11 = 1
    12 = 2
    13 = 2
    14 = 2
    15 = 3
    16 = 2
    17 = 1
    18 = 1
    19 = 1
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 2
    25 = 1
    26 = 1
    27 = 2
    28 = 0
    29 = 1
    30 = 1
    31 = 1
    32 = 1
    33 = 1
    34 = 3
    35 = 1
    36 = 0
    37 = 3
    38 = 0
    39 = 3
    40 = 2
    41 = 2
    42 = 2
    43 = 3
    44 = 2
    45 = 1
    46 = 3
    47 = 0
    48 = 1
    49 = 0
    50 = 0
    51 = 1
    52 = 2
    53 = 2
    54 = 2
    55 = 1
    56 = 1
    57 = 1
    58 = 0
    59 = 0
    60 = 1
    61 = 3
    62 = 1
    63 = 2
    64 = 2
    65 = 3
    66 = 0
    67 = 1
    68 = 0
    69 = 2
    70 = 0
    71 = 2
    72 = 1
    73 = 2
    74 = 2
    75 = 1
    76 = 1
    77 = 0
    78 = 1
    79 = 1
    80 = 1
    81 = 1
    82 = 1
    83 = 3
    84 = 5
    85 = 1
    87 = 3
    88 = 2
    89 = 2
    91 = 2
    92 = 1
    93 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
